$wb = $excel.ActiveWorkbook

$wsSrc = $wb.Worksheets.Item("On-Site")
$wsDst = $wb.Worksheets.Item("Replacement")

# Copy the maintenance rows (2-5) from On-Site down to the Replacement sheet,
# mirroring the same layout (columns A:V).
$srcRange = $wsSrc.Range("A2:V5")
$srcRange.Copy()
$dstRange = $wsDst.Range("A2")
$dstRange.PasteSpecial()

# Widen column A on Replacement to fit the new entries.
$wsDst.Columns.Item(1).ColumnWidth = 16.81640625

# Update selections on both sheets.
$wsSrc.Select()
$wsSrc.Range("A2:XFD5").Select()

$wsDst.Select()
$wsDst.Range("A2:XFD5").Select()
